$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("condensed")

$values = @{
    "B2"  = 294714441.5634915
    "C2"  = 148607180.505959
    "B3"  = 119773340.2477307
    "C3"  = 59828713.9764261
    "B4"  = 97333301.33751126
    "C4"  = 45984910.43738856
    "B5"  = 114174629.4046326
    "C5"  = 53574611.39772989
    "B6"  = 142742068.5514288
    "C6"  = 64388111.45899402
    "B7"  = 138797161.9538306
    "C7"  = 58332255.54209531
    "B8"  = 115297238.2501935
    "C8"  = 44508339.75771271
    "B9"  = 139481601.8581507
    "C9"  = 61463725.2981982
    "B10" = 178787506.9410357
    "C10" = 88516124.21678902
    "B11" = 198295855.3685375
    "C11" = 102189588.0096093
    "B12" = 189240484.9584804
    "C12" = 100792107.9889705
    "B13" = 176997595.0678757
    "C13" = 94690949.00786746
    "B14" = 167744663.6431007
    "C14" = 89636205.40607382
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
